# Update the "想去人数" (want-to-go count) values in column F on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets to reflect the
# latest generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 78
$ws1.Range("F7").Value = 7640
$ws1.Range("F8").Value = 85
$ws1.Range("F13").Value = 438
$ws1.Range("F14").Value = 162
$ws1.Range("F16").Value = 429
$ws1.Range("F17").Value = 59
$ws1.Range("F18").Value = 58
$ws1.Range("F20").Value = 5495
$ws1.Range("F21").Value = 145
$ws1.Range("F22").Value = 203
$ws1.Range("F23").Value = 923
$ws1.Range("F25").Value = 301

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 78
$ws4.Range("F7").Value = 7640
$ws4.Range("F8").Value = 85
$ws4.Range("F13").Value = 438
$ws4.Range("F14").Value = 162
$ws4.Range("F16").Value = 429
$ws4.Range("F17").Value = 59
$ws4.Range("F18").Value = 58
$ws4.Range("F21").Value = 5495
$ws4.Range("F23").Value = 145
$ws4.Range("F24").Value = 203
$ws4.Range("F25").Value = 923
$ws4.Range("F27").Value = 301
